$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain text (matches original inlineStr formatting)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "96.739.97"
$ws.Range("E2").Value = "  -1.99%  "

$ws.Range("D3").Value = "3.295.73"
$ws.Range("E3").Value = "  -5.20%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "245.56"
$ws.Range("E5").Value = "  -5.91%  "

$ws.Range("D6").Value = "646.39"
$ws.Range("E6").Value = "  -3.95%  "

$ws.Range("D7").Value = "1.35"
$ws.Range("E7").Value = "  -13.79%  "

$ws.Range("D8").Value = "0.408"
$ws.Range("E8").Value = "  -11.53%  "

$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.03%  "

$ws.Range("D10").Value = "0.979"
$ws.Range("E10").Value = "  -12.24%  "

$ws.Range("D11").Value = "3.287.75"
$ws.Range("E11").Value = "  -5.48%  "

$ws.Range("D12").Value = "0.203"
$ws.Range("E12").Value = "  -8.98%  "

$ws.Range("D13").Value = "39.47"
$ws.Range("E13").Value = "  -7.53%  "

$ws.Range("D14").Value = "96.501.30"
$ws.Range("E14").Value = "  -1.56%  "

$ws.Range("D15").Value = "5.92"
$ws.Range("E15").Value = "  -3.07%  "

$ws.Range("D16").Value = "0.0000247"
$ws.Range("E16").Value = "  -9.70%  "

$ws.Range("D17").Value = "3.902.86"
$ws.Range("E17").Value = "  -5.31%  "

$ws.Range("D18").Value = "8.58"
$ws.Range("E18").Value = "  +4.31%  "

$ws.Range("D19").Value = "3.293.18"
$ws.Range("E19").Value = "  -5.11%  "

$ws.Range("D20").Value = "16.47"
$ws.Range("E20").Value = "  -6.12%  "

$ws.Range("D21").Value = "0.507"
$ws.Range("E21").Value = "  +10.32%  "

$ws.Range("D22").Value = "10.31"
$ws.Range("E22").Value = "  -4.91%  "

$ws.Range("D23").Value = "486.41"
$ws.Range("E23").Value = "  -10.08%  "

$ws.Range("D24").Value = "3.23"
$ws.Range("E24").Value = "  -11.17%  "

$ws.Range("D25").Value = "0.0000194"
$ws.Range("E25").Value = "  -12.47%  "

$ws.Range("D26").Value = "6.25"
$ws.Range("E26").Value = "  -2.17%  "

$ws.Range("D27").Value = "91.98"
$ws.Range("E27").Value = "  -10.75%  "

$ws.Range("D28").Value = "11.81"
$ws.Range("E28").Value = "  -8.56%  "

$ws.Range("D29").Value = "3.481.34"
$ws.Range("E29").Value = "  -4.53%  "

$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").Value = "0.991"
$ws.Range("E30").Value = "  -0.77%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.141"
$ws.Range("E31").Value = "  -7.20%  "

$ws.Range("D32").Value = "10.64"
$ws.Range("E32").Value = "  -7.91%  "

$ws.Range("D33").Value = "0.184"
$ws.Range("E33").Value = "  -9.20%  "

$ws.Range("D34").Value = "2.42"
$ws.Range("E34").Value = "  +8.63%  "

$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("D36").Value = "0.535"
$ws.Range("E36").Value = "  -7.21%  "

$ws.Range("D37").Value = "27.70"
$ws.Range("E37").Value = "  -10.09%  "

$ws.Range("D38").Value = "1.43"
$ws.Range("E38").Value = "  -0.45%  "

$ws.Range("B39").Value = "USDe"
$ws.Range("C39").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").Value = "7.37"
$ws.Range("E40").Value = "  -7.72%  "

$ws.Range("D41").Value = "0.147"
$ws.Range("E41").Value = "  -8.65%  "

$ws.Range("D42").Value = "494.81"
$ws.Range("E42").Value = "  -8.51%  "

$ws.Range("D43").Value = "24.55"
$ws.Range("E43").Value = "  -0.91%  "

$ws.Range("D44").Value = "3.66"
$ws.Range("E44").Value = "  -3.45%  "

$ws.Range("D45").Value = "0.811"
$ws.Range("E45").Value = "  -7.07%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0403"
$ws.Range("E46").Value = "  -8.34%  "

$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").Value = "8.33"
$ws.Range("E47").Value = "  +0.96%  "

$ws.Range("D48").Value = "1.60"
$ws.Range("E48").Value = "  +1.07%  "

$ws.Range("D49").Value = "5.33"
$ws.Range("E49").Value = "  +0.89%  "

$ws.Range("D50").Value = "52.19"
$ws.Range("E50").Value = "  +1.60%  "

$ws.Range("D51").Value = "3.08"
$ws.Range("E51").Value = "  -12.49%  "
